$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Jonas Valanciunas", "C", "Sacramento Kings"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Matas Buzelis", "SF,PF", "Chicago Bulls"),
    @("Jimmy Butler III", "SF,PF", "Golden State Warriors"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
